$wb = $excel.ActiveWorkbook

# The new sheet mirrors Hoja3's layout/formatting, so create it by copying
# Hoja3 and dropping it after the last sheet - this keeps the same XML
# namespace declarations (mc:Ignorable="x14ac", etc.) that a hand-authored
# worksheet part carries.
$srcSheet = $wb.Worksheets.Item(3)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$srcSheet.Copy($null, $lastSheet)

$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "Hoja4"
$ws.Range("A1").Value = "otherside"

# Move the selection to A2 (as Excel does after committing a value in A1)
# and make the new sheet the active one.
[void]$ws.Range("A2").Select()
$ws.Activate()
